$wb = $excel.ActiveWorkbook

# Work on the ParserWriter sheet
$ws = $wb.Worksheets.Item("ParserWriter")

# Row 22 was already a blank (but styled) spacer row, so it simply gets
# filled in with the new "run management record" entry in place; a new
# row is inserted above the next table's header (old row 23, "PyST.feflow
# Package") to push that table down by one row and make room.
$ws.Rows.Item(23).Insert()

# Fill the (already existing, previously blank) row 22 with the new entry
$ws.Cells.Item(22, 1).Value = "run management record "
$ws.Cells.Item(22, 2).Value = "rmr"
$ws.Cells.Item(22, 3).Value = 1

# Update the selection on this sheet and make it the active tab
$ws.Range("E21").Select()
$ws.Activate()
